$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata" (sheet1)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# B7 (Experimental value) goes from empty to the text "true".
# A plain assignment of the literal text "true"/"false" is auto-coerced to a
# Boolean by Excel, so instead we enter it with a leading quote-prefix (forces
# text) and then copy the number format from a neighboring "normal" data cell
# back onto it so the cell ends up with the same style as every other data
# row (instead of picking up the quote-prefixed style).
$b7 = $meta.Cells.Item(7, 2)
$b7.Value = "'true"
$meta.Cells.Item(6, 2).Copy()
$b7.PasteSpecial(-4122)   # xlPasteFormats

# B8: Date
$meta.Cells.Item(8, 2).Value = "2024-12-27T22:28:35+00:00"

# B12: Description
$meta.Cells.Item(12, 2).Value = "Value set for measurement quality indicators"

# ---------------------------------------------------------------------------
# Sheet "Include #0" (sheet2)
# ---------------------------------------------------------------------------
$inc = $wb.Worksheets.Item("Include #0")

# Original rows:
#  1 Concept     | Description
#  2 723510000   | High quality
#  3 723511001   | Moderate quality
#  4 723512008   | Low quality
#  5 723513003   | Uncertain quality
#  6 (empty)     | (empty)
#  7 System URI  | http://snomed.info/sct
#
# Target rows:
#  1 Codes       | (none)
#  2 All codes   | (none)
#  3 (empty)     | (empty)
#  4 System URI  | https://github.com/RicardoLSantos/shorthand/CodeSystem/measurement-quality-cs
#
# Remove rows 3-5 so the pre-existing blank row (old row 6) slides up to row 3
# and the System URI row (old row 7) slides up to row 4 - this preserves the
# blank row's existing (shared-string) empty-text cells exactly.
$inc.Rows.Item(5).Delete()
$inc.Rows.Item(4).Delete()
$inc.Rows.Item(3).Delete()

# Row 1: Concept -> Codes, and clear column B entirely (no cell at all).
$inc.Cells.Item(1, 1).Value = "Codes"
$inc.Cells.Item(1, 2).Clear()

# Row 2: first concept code -> All codes, and clear column B entirely.
$inc.Cells.Item(2, 1).Value = "All codes"
$inc.Cells.Item(2, 2).Clear()

# Row 3 is already the untouched, pre-existing blank row - nothing to do.

# Row 4 (previously row 7): keep "System URI" label, update the URL value.
$inc.Cells.Item(4, 1).Value = "System URI"
$inc.Cells.Item(4, 2).Value = "https://github.com/RicardoLSantos/shorthand/CodeSystem/measurement-quality-cs"
